# Update row 12 (Fiscal Year 2081/82) with the complete fiscal year figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 13.25
$ws.Range("C12").Value = 81.8
$ws.Range("D12").Value = 6
$ws.Range("E12").Value = 19.24
$ws.Range("F12").Value = -37.7
$ws.Range("G12").Value = 52.47
$ws.Range("H12").Value = -5.02
